$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A56").NumberFormat = "@"
$ws.Range("A56").Value = "2025-10-10"
$ws.Range("A56").ClearFormats()
$ws.Range("B56").Value = 54.31000137329102
$ws.Range("C56").Value = 678.9500122070312
$ws.Range("D56").Value = 348.2999877929688
